$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.688.87'
$ws.Range('E2').Value = '  -0.61%  '

$ws.Range('D3').Value = '2.455.24'
$ws.Range('E3').Value = '  -0.69%  '

$ws.Range('E4').Value = '  +0.01%  '

$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '570.76'
$c.ClearFormats()

$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '146.19'
$c.ClearFormats()
$ws.Range('E6').Value = '  -0.41%  '

$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('E8').Value = '  -1.90%  '

$ws.Range('E9').Value = '  -1.42%  '

$ws.Range('E10').Value = '  -0.12%  '

$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '5.17'
$c.ClearFormats()
$ws.Range('E11').Value = '  -2.13%  '

$ws.Range('E12').Value = '  -1.91%  '

$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '28.52'
$c.ClearFormats()
$ws.Range('E13').Value = '  -0.84%  '

$ws.Range('E14').Value = '  -3.68%  '

$ws.Range('D15').Value = '2.900.66'
$ws.Range('E15').Value = '  -0.66%  '

$ws.Range('D16').Value = '62.589.77'
$ws.Range('E16').Value = '  -0.56%  '

$ws.Range('D17').Value = '2.454.09'
$ws.Range('E17').Value = '  -0.84%  '

$ws.Range('E19').Value = '  -3.07%  '

$ws.Range('B20').Value = 'SuiNetwork'
$ws.Range('C20').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '2.23'
$c.ClearFormats()
$ws.Range('E20').Value = '  -0.67%  '

$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '321.00'
$c.ClearFormats()
$ws.Range('E21').Value = '  -2.58%  '

$ws.Range('E22').Value = '  -0.26%  '

$ws.Range('E23').Value = '  -0.03%  '

$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '9.93'
$c.ClearFormats()
$ws.Range('E24').Value = '  +3.32%  '

$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '64.68'
$c.ClearFormats()
$ws.Range('E25').Value = '  -2.51%  '

$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '646.14'
$c.ClearFormats()
$ws.Range('E26').Value = '  -3.38%  '

$ws.Range('D27').Value = '2.575.70'
$ws.Range('E27').Value = '  -0.67%  '

$ws.Range('D28').Value = '0.0₃0951'
$ws.Range('E28').Value = '  -3.99%  '

$ws.Range('E29').Value = '  -0.19%  '

$ws.Range('E30').Value = '  -3.70%  '

$ws.Range('E31').Value = '  -2.92%  '

$ws.Range('E32').Value = '  -3.49%  '

$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '0.133'
$c.ClearFormats()
$ws.Range('E33').Value = '  -0.42%  '

$ws.Range('E34').Value = '  -0.02%  '

$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '1.48'
$c.ClearFormats()
$ws.Range('E35').Value = '  -4.13%  '

$ws.Range('E36').Value = '  -3.35%  '

$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '150.45'
$c.ClearFormats()
$ws.Range('E37').Value = '  -0.01%  '

$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '18.54'
$c.ClearFormats()
$ws.Range('E38').Value = '  -1.32%  '

$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '0.363'
$c.ClearFormats()
$ws.Range('E39').Value = '  -2.43%  '

$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '5.31'
$c.ClearFormats()
$ws.Range('E40').Value = '  -3.16%  '

$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '2.63'
$c.ClearFormats()
$ws.Range('E41').Value = '  -3.90%  '

$ws.Range('E42').Value = '  -3.73%  '

$ws.Range('D43').Value = '0.0₆0310'
$ws.Range('E43').Value = '  +1.36%  '

$ws.Range('E44').Value = '  +0.62%  '

$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '152.66'
$c.ClearFormats()
$ws.Range('E45').Value = '  -0.04%  '

$ws.Range('E46').Value = '  +1.79%  '

$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '3.53'
$c.ClearFormats()
$ws.Range('E47').Value = '  -2.02%  '

$ws.Range('E48').Value = '  -0.58%  '

$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '19.93'
$c.ClearFormats()
$ws.Range('E49').Value = '  -3.45%  '

$ws.Range('E50').Value = '  -1.60%  '

$ws.Range('E51').Value = '  -1.88%  '
